# Update EUR->ARS rate: add new reading row 2025-10-10T21:21:54Z
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 70

# Column A holds a plain text date (e.g. "2025-10-10"), not a real Excel
# date. Force text formatting on that cell first so Excel does not
# auto-convert the string into a date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-10"
$ws.Cells.Item($row, 2).Value = "21:21:54"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,750.2781"
